# Updates cryptos list values (price & volume) per commit "Updated cryptos list on Sat Dec 16 13:21:48 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.504.92"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "'2.252.84"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'246.29"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "'76.18"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("E10").Value = "  +7.28%  "
$ws.Range("D11").Value = "'0.0952"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "'7.32"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "'2.595.36"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "'14.66"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "'0.858"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'2.270.33"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'42.334.72"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "'0.0000101"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'72.43"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").Value = "'2.23"
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("D23").Value = "'232.01"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'9.10"
$ws.Range("E24").Value = "  +33.93%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'11.44"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("D27").Value = "'3.61"
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'168.70"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "'20.74"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "'31.11"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "'5.32"
$ws.Range("E35").Value = "  +12.44%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "'4.53"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "'0.0319"
$ws.Range("E38").Value = "  +6.52%  "
$ws.Range("D39").Value = "'13.69"
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("D40").Value = "'2.19"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").Value = "'5.84"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Value = "'63.26"
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("D43").Value = "'0.203"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").Value = "'108.51"
$ws.Range("E44").Value = "  -8.28%  "
$ws.Range("D45").Value = "'8.80"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "'2.34"
$ws.Range("E50").Value = "  +4.92%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "'4.14"
$ws.Range("E51").Value = "  -6.27%  "
